$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data (player -> position/team) is unchanged; only the
# row order of the player table (A2:C17) has been resorted.
$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Mike Conley", "PG", "Minnesota Timberwolves"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Brandon Boston Jr.", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
